$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = -0.92012669354727072
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value = -0.3863929271143855

$ws.Range("B3").Value = -1.3172335172622707
$ws.Range("C3").Value = 0.67302120142693267
$ws.Range("D3").Value = -0.68186776113077985
$ws.Range("E3").Value = 2.1324821486706291

$ws.Range("B1:E3").Select()
